# Update cryptos list figures (prices and 1h volume %) per data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.647.45'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.112.77'
$ws.Range('E3').Value = '  +9.81%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''254.61'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('E6').Value = '  -4.74%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '''47.14'
$ws.Range('E8').Value = '  +5.88%  '
$ws.Range('D9').Value = '''60.61'
$ws.Range('E9').Value = '  +3.40%  '
$ws.Range('D10').Value = '''0.373'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').Value = '''0.0745'
$ws.Range('E11').Value = '  -2.58%  '
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = '2.424.59'
$ws.Range('E13').Value = '  +10.08%  '
$ws.Range('D14').Value = '''14.25'
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('D15').Value = '''0.834'
$ws.Range('E15').Value = '  +3.72%  '
$ws.Range('D16').Value = '2.112.61'
$ws.Range('E16').Value = '  +9.88%  '
$ws.Range('D17').Value = '''5.11'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').Value = '36.587.78'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').Value = '''73.66'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('D20').Value = '0.0₃0835'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').Value = '''13.22'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').Value = '''240.95'
$ws.Range('E22').Value = '  -4.51%  '
$ws.Range('D23').Value = '''5.19'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  -7.18%  '
$ws.Range('D26').Value = '''171.83'
$ws.Range('E26').Value = '  +2.00%  '
$ws.Range('D27').Value = '''21.72'
$ws.Range('E27').Value = '  +15.20%  '
$ws.Range('D28').Value = '''9.18'
$ws.Range('E28').Value = '  +3.68%  '
$ws.Range('D29').Value = '''2.01'
$ws.Range('E29').Value = '  -8.84%  '
$ws.Range('D30').Value = '''29.26'
$ws.Range('E30').Value = '  +64.05%  '
$ws.Range('D31').Value = '''0.123'
$ws.Range('E31').Value = '  -4.74%  '
$ws.Range('D32').Value = '''4.48'
$ws.Range('E32').Value = '  -1.77%  '
$ws.Range('D33').Value = '''0.0950'
$ws.Range('E33').Value = '  +9.48%  '
$ws.Range('D34').Value = '''0.0600'
$ws.Range('E34').Value = '  -3.41%  '
$ws.Range('D35').Value = '''0.953'
$ws.Range('E35').Value = '  +5.81%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.89'
$ws.Range('E36').Value = '  -4.26%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '''2.34'
$ws.Range('E37').Value = '  +15.71%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = '''4.14'
$ws.Range('E39').Value = '  -5.39%  '
$ws.Range('D40').Value = '''1.34'
$ws.Range('E40').Value = '  -12.23%  '
$ws.Range('D41').Value = '''1.19'
$ws.Range('E41').Value = '  +6.74%  '
$ws.Range('D42').Value = '''0.0225'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('D43').Value = '''98.76'
$ws.Range('E43').Value = '  -7.38%  '
$ws.Range('E44').Value = '  +8.50%  '
$ws.Range('D45').Value = '''16.02'
$ws.Range('E45').Value = '  -7.81%  '
$ws.Range('D46').Value = '1.350.86'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('D47').Value = '''7.21'
$ws.Range('E47').Value = '  +11.52%  '
$ws.Range('D48').Value = '''0.0841'
$ws.Range('E48').Value = '  +3.10%  '
$ws.Range('D49').Value = '2.295.44'
$ws.Range('E49').Value = '  +9.30%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '''2.83'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '''2.28'
$ws.Range('E51').Value = '  -4.62%  '
